$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Vegas line updates (home favourite spread / over-under totals) for NE (row 2)
# and GB (row 5), plus swapping the row-4 pick from MIN to ATL with its own
# updated spread. Everything else (E/F/G/H/I/J columns, the B9:C12 summary
# strings, etc.) is formula-driven and recalculates automatically.
$ws.Range("C2").Value = -4
$ws.Range("D2").Value = 41

$ws.Range("B4").Value = "ATL"
$ws.Range("C4").Value = -5

$ws.Range("C5").Value = -3
$ws.Range("D5").Value = 40
